$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting the existing rows 26-45 down to 27-46.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the week's new data entry.
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 45202
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 100112013
$ws.Range("G26").Value = "Alcachofa"
$ws.Range("H26").Value = "Española"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 12000
$ws.Range("N26").Value = "$/caja 30 unidades"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 400
$ws.Range("Q26").Value = 30
$ws.Range("R26").Value = "Hortaliza"
